$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "zeroshot huang combined with own"
$ws.Range("D13").Value = 3000
$ws.Range("E13").Value = 200
$ws.Range("F13").Value = 5
$ws.Range("I13").Value = "auto"
$ws.Range("J13").Value = "yes"

$ws.Range("L13").Select()
